# Generate Report for Handback
# This script reproduces the "handback" report-generation edit:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - For each localized-language sheet (zh-cn, de-de), the "Latest Target File" (I)
#    and "Latest Handback File" (J) columns get populated, I becomes a hyperlink to the
#    source markdown file (like column A), and the "Latest Handback DateTime" (K) is stamped.
#  - Column widths are widened for the affected columns.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 1) Status column text: "Ready for handoff" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2) Populate "Latest Target File" (I) / "Latest Handback File" (J) /
#    "Latest Handback DateTime" (K) for zh-cn and de-de sheets, and add
#    hyperlinks on column I (same targets as column A's hyperlinks).
# ---------------------------------------------------------------------------

$url27f5425e = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/53eceb53c78fec6383f8a3aeda9ec7ff8dd81d65/e2e/27f5425e-7377-4959-9110-1f54699a9831.md"
$url2f7db598 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/53eceb53c78fec6383f8a3aeda9ec7ff8dd81d65/e2e/2f7db598-7b81-4391-a6ab-0a7ed8fae673.md"

$disp27f5425e = "27f5425e-7377-4959-9110-1f54699a9831.md"
$disp2f7db598 = "2f7db598-7b81-4391-a6ab-0a7ed8fae673.md"

# -- zh-cn sheet --
$wsZhCn.Range("J2").Value = "27f5425e-7377-4959-9110-1f54699a9831.518a67538865fa12b14fb87deb7d674cd9a0135d.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-25 09:06:34"

$wsZhCn.Range("J3").Value = "2f7db598-7b81-4391-a6ab-0a7ed8fae673.99e9a99aeeb7aca0e52c690894766408abdc79e1.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-25 09:06:34"

# Rebuild hyperlinks collection on zh-cn: A2, I2, A3, I3 (this reproduces the
# rId2/rId3/rId4/rId5 relationship ordering used by the handback report tool).
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $url27f5425e, "", "", $disp27f5425e) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $url27f5425e, "", "", $disp27f5425e) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $url2f7db598, "", "", $disp2f7db598) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $url2f7db598, "", "", $disp2f7db598) | Out-Null

# -- de-de sheet --
$wsDeDe.Range("J2").Value = "27f5425e-7377-4959-9110-1f54699a9831.518a67538865fa12b14fb87deb7d674cd9a0135d.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-25 09:06:41"

$wsDeDe.Range("J3").Value = "2f7db598-7b81-4391-a6ab-0a7ed8fae673.99e9a99aeeb7aca0e52c690894766408abdc79e1.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-25 09:06:41"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $url27f5425e, "", "", $disp27f5425e) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $url27f5425e, "", "", $disp27f5425e) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $url2f7db598, "", "", $disp2f7db598) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $url2f7db598, "", "", $disp2f7db598) | Out-Null

# ---------------------------------------------------------------------------
# 3) Column width adjustments
# ---------------------------------------------------------------------------
# Overview columns E (zh-cn) and F (de-de) grow to fit the new status text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

# zh-cn / de-de: column C (Status) grows; columns I/J (Latest Target File /
# Latest Handback File) grow to the standard 40-character width.
$wsZhCn.Columns.Item(3).ColumnWidth = 29.1
$wsZhCn.Columns.Item(9).ColumnWidth = 39.1667
$wsZhCn.Columns.Item(10).ColumnWidth = 39.1667

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1
$wsDeDe.Columns.Item(9).ColumnWidth = 39.1667
$wsDeDe.Columns.Item(10).ColumnWidth = 39.1667
